# Edit LOQ4231.xlsx: remove the "Docentes responsáveis" row, remove the long
# "Programa" paragraph row, and shuffle a few remaining labels/values so the
# sheet ends up with 21 rows (A1:C21) instead of 22 (A1:C22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (old row 22, the "Bibliografia" reference text) -
# everything below row 21 shifts up by one, giving the new 21-row sheet.
$ws.Rows.Item(22).Delete()

# --- Row 10: "Objetivos:" row now shows the teacher's name instead of the
#     old free-text objective paragraph.
$ws.Range("B10").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C10").Value = "5840671 - Francisco José Moreira Chaves"

# --- Row 13 (previously just had B/C holding the teacher's name): becomes
#     "Programa resumido:" / "Semestral".
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14: becomes "Short syllabus:" with no B/C value.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# --- Row 15: becomes "Programa:" with a date value in B/C.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

# --- Row 16: becomes "Syllabus:" with no B/C value.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# --- Row 17: becomes "Avaliação:" (no B/C).
$ws.Range("A17").Value = "Avaliação:"

# --- Row 18: becomes "Método:" with the teacher's name repeated in B/C.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C18").Value = "5840671 - Francisco José Moreira Chaves"

# --- Row 19: label becomes "Critério:" (B/C keep the "Aulas Expositivas..." text).
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label becomes "Norma de recuperação:" (B/C keep the "MF = ..." text).
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label becomes "Bibliografia:" (B/C keep the "NF = ..." text).
$ws.Range("A21").Value = "Bibliografia:"

# --- Row-height cleanup ---------------------------------------------------
# Rows 10, 11 and 17 lose their explicit 60/120pt custom height and go back
# to the default row height.
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(17).AutoFit()

# Rows 13, 14, 18, 19, 20 get an explicit 60pt height; rows 15, 16, 21 get
# an explicit 120pt height.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
